$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Status column: "In Translation" -> "Ready for handoff"
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# Overview "Latest Handoff Date"
$wsOverview.Range("D2").Value = "2016-28-17 14:28:17"

# zh-cn "Latest Handoff Datetime"
$wsZhCn.Range("E2").Value = "2016-03-17 14:28:14"

# de-de "Latest Handoff Datetime"
$wsDeDe.Range("E2").Value = "2016-03-17 14:28:17"
